$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0,"falling",-2.246021509170532,2.530338048934937,1.034301340579987,-0.0774271711707115,0.1702786833047866,-0.0888808965682983),
  @(100,"falling",-2.397771739959717,2.518170547485352,0.915059447288513,-0.047036625444889,0.0433714315295219,-0.0064140851609408),
  @(200,"falling",-2.324807786941528,2.547794580459595,0.7709671020507812,0.0061086523346602,-0.0146607663482427,-0.0610865242779254),
  @(300,"falling",-2.206723213195801,2.590793609619141,0.8417039632797245,0.0401643887162208,0.0106901414692401,0.0001527163112768),
  @(400,"falling",-2.01174955368042,2.60823769569397,1.090024280548096,0.0146607663482427,0.0265726372599601,0.0464257597923278),
  @(500,"falling",-2.355203628540039,2.568441867828369,0.8857938051223755,-0.0455094613134861,0.0074830991216003,0.0467311926186084),
  @(600,"falling",-2.247145366668701,2.438617801666259,0.8904013574123383,-0.0215329993516206,0.0041233403608202,0.0065668015740811),
  @(700,"falling",-2.348918724060059,2.538596057891846,1.142265951633454,0.0096211275085806,-0.0128281703218817,0.0253509078174829),
  @(800,"falling",-2.434140348434448,2.598772668838501,1.196871364116669,0.024892758578062,-0.0320704244077205,0.0131336031481623),
  @(900,"falling",-2.469630908966064,2.608709144592285,0.9890874266624448,0.0187841057777404,-0.0146607663482427,-0.0226020142436027),
  @(1000,"falling",-2.473813533782959,2.590450048446655,0.9850140511989596,-0.028557950630784,-0.0305432621389627,-0.0201585534960031),
  @(1100,"falling",-2.382417774200439,2.536800479888916,1.009570789337157,-0.0192422550171613,0.0343611687421798,-0.0004581489483825),
  @(1200,"falling",-2.251543283462525,2.502760457992554,0.9712224066257489,0.009010262787342,0.0630718395113945,-0.0142026171088218),
  @(1300,"falling",-2.286794948577881,2.536056137084961,1.293239164352417,0.0442877300083637,0.0134390350431203,-0.001527163083665),
  @(1400,"falling",-2.406369590759277,2.443658208847046,1.501735496520997,-0.0262672062963247,-0.027030786499381,-0.0058032199740409),
  @(1500,"falling",-2.325653553009033,2.315335750579834,2.114412546157837,-0.0975857228040695,0.0615446716547012,-0.0801760628819465),
  @(1600,"falling",-1.993824815750122,2.412744569778442,2.562826609611511,-0.0795651972293853,0.1488984078168869,-0.1036943718791008),
  @(1700,"falling",-2.099692916870118,2.386767387390137,2.678394412994385,0.0178678091615438,0.0836885422468185,0.0500909499824047),
  @(1800,"falling",-2.539734458923342,2.251132488250732,2.302260518074033,0.0439822971820831,0.1259909570217132,-0.0801760628819465),
  @(1900,"falling",-2.798751735687256,1.147012233734121,2.277631902694703,0.1090394482016563,0.9398161768913268,-0.2848159074783325),
  @(2000,"falling",-3.073579788208008,-0.642811775207524,3.430169939994816,0.0601702257990837,0.5073235630989075,0.0161879286170005),
  @(2100,"falling",-4.291536808013922,-2.572477817535406,5.254911422729497,-0.3770565688610077,0.2907718420028686,-1.639104127883911),
  @(2200,"falling",-6.406326150894171,-2.793146014213552,5.593617057800286,-1.575268745422363,0.9181304574012756,-3.197726726531982),
  @(2300,"falling",0.9536849975587334,-0.1826061248779267,2.677957320213301,0.5009095072746277,3.925878286361694,0.5390886068344116),
  @(2400,"falling",32.34990999698625,0.490568253397939,-0.6328019738197053,2.016924381256104,-0.4338670372962951,-2.551736831665039),
  @(2500,"falling",-1.651923894882202,-0.1518822014331817,6.22404146194458,-1.234558701515198,2.077705383300781,-0.801913321018219),
  @(2600,"falling",-4.871152138710036,1.64225258529187,4.566823816299431,-0.8013024926185608,-1.52746856212616,-0.2012801021337509),
  @(2700,"falling",-6.492228031158454,3.64010591506959,1.834508514404282,-0.7525860071182251,-3.772551059722901,1.262963891029358),
  @(2800,"falling",-1.539169025421104,1.307969903945902,4.571132516860986,-0.4665483236312866,0.3252857327461242,0.3428481221199035),
  @(2900,"falling",-1.01956577301027,2.104999399185199,2.938691687583895,-0.1076650023460388,-0.3216205537319183,-0.2797762751579284)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}
